# Update PoFDCtAE worksheet: replace several formula-driven cells with
# hardcoded value 1, which in turn changes the computed value of the
# dependent "1-x" formulas in column R.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoFDCtAE")

$ws.Range("C3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("N14").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("T20").Value = 1

$excel.Calculate()
